$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update StudentID values
$ws.Range("B2").Value = 2073777
$ws.Range("B3").Value = 2071398

# Update the active selection to D12
$ws.Range("D12").Select()
